# Update "想去人数" (F column) figures on both the "展览" sheet and the
# "全部类型" sheet, which mirrors the same rows.
$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1617
    3  = 9003
    8  = 180
    11 = 3842
    22 = 244
    24 = 2651
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
